$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = 'user shares projects in publicly available repository - devin says "project sharing privileges"'
$ws.Range("A15").Select()
